# COOLMATE-005.xlsx — "executed Testsuite COOLMATE005 successfully"
#
# Two substantive edits on the single worksheet:
#   1. Fix a typo in the Test Case name for scenario #3: "TC-003-Delete
#      address" -> "TC003-Delete address" (stray hyphen removed).
#   2. Leave the sheet's cursor/viewport the way it was when the run
#      finished: cell E7 selected, with the view scrolled so column B is
#      the left-most visible column (topLeftCell = B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the Test Case label in row 4 (column E, "Test Case").
$ws.Range("E4").Value = "TC003-Delete address"

# 2) Reflect the final selection/scroll position left by the test run.
$ws.Range("E7").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
